$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the shadow-rate re-estimates for the existing history (rows 83-128, 147-152)
$ws.Range("C83").Value = 1.0095613250623581
$ws.Range("C84").Value = 1.432656869497051
$ws.Range("C85").Value = 1.9491960683761356
$ws.Range("C86").Value = 2.4691396904516605
$ws.Range("C87").Value = 2.9424669894994926
$ws.Range("C88").Value = 3.4591638057257734
$ws.Range("C89").Value = 3.9792204408714893
$ws.Range("C90").Value = 4.4559627369300747
$ws.Range("C91").Value = 4.9060510801289148
$ws.Range("C92").Value = 5.2461475276541103
$ws.Range("C93").Value = 5.2462486955683785
$ws.Range("C94").Value = 5.2563500745136826
$ws.Range("C95").Value = 5.2497827934618213
$ws.Range("C96").Value = 5.0732119369917728
$ws.Range("C97").Value = 4.4966360198162914
$ws.Range("C98").Value = 3.1767200416648711
$ws.Range("C99").Value = 2.0867965552669343
$ws.Range("C100").Value = 1.9402002012724795
$ws.Range("C101").Value = 0.50805694283571068
$ws.Range("C102").Value = 1.6143332037749536
$ws.Range("C103").Value = 0.20491834297065026
$ws.Range("C104").Value = -0.71108952572079831
$ws.Range("C105").Value = -0.76912085303370015
$ws.Range("C106").Value = -0.60252808630845145
$ws.Range("C107").Value = -2.0634282890471534
$ws.Range("C108").Value = -1.8886723178182963
$ws.Range("C109").Value = -2.5300849646885615
$ws.Range("C110").Value = -1.9808665525398617
$ws.Range("C111").Value = -1.5854631537771002
$ws.Range("C112").Value = -2.8607062990615661
$ws.Range("C113").Value = -2.4217347958699897
$ws.Range("C114").Value = -3.313265728242154
$ws.Range("C115").Value = -3.0897822715132306
$ws.Range("C116").Value = -2.6436980305593272
$ws.Range("C117").Value = -3.857697797500248
$ws.Range("C118").Value = -2.3271712868010663
$ws.Range("C119").Value = -1.6049835861445816
$ws.Range("C120").Value = -1.211003789002163
$ws.Range("C121").Value = -1.2960217529530471
$ws.Range("C122").Value = -1.3459005823200587
$ws.Range("C123").Value = -1.277261379068606
$ws.Range("C124").Value = -0.8669078890010784
$ws.Range("C125").Value = -0.40434927102547258
$ws.Range("C126").Value = 0.131373396207346
$ws.Range("C127").Value = 0.017823339012745798
$ws.Range("C128").Value = 0.04899030447687025
$ws.Range("C147").Value = 7.7328677810758473
$ws.Range("C148").Value = -6.942753656964995
$ws.Range("C149").Value = -4.4147876959535148
$ws.Range("C150").Value = -3.615120344168421
$ws.Range("C151").Value = -3.1359489682822206
$ws.Range("C152").Value = -1.6439898531997899

# Append the newest observation (2021 Q4) to the bottom of the series
$ws.Range("A153").Value = 2021.75
$ws.Range("B153").Value = 0
$ws.Range("C153").Value = -0.65474866359127581

